# Apply score/odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Cells.Item(3, 8).Value  = 5.2    # H3: 5.1 -> 5.2
$ws.Cells.Item(3, 10).Value = 4.6    # J3: 4.7 -> 4.6

# Row 4
$ws.Cells.Item(4, 7).Value  = 3.1    # G4: 3.35 -> 3.1
$ws.Cells.Item(4, 8).Value  = 2.44   # H4: 2.42 -> 2.44
$ws.Cells.Item(4, 16).Value = 2.18   # P4: 2.2 -> 2.18
$ws.Cells.Item(4, 17).Value = 1.7    # Q4: 1.69 -> 1.7

# Row 5
$ws.Cells.Item(5, 6).Value  = 3.45   # F5: 3.4 -> 3.45
$ws.Cells.Item(5, 11).Value = 4.8    # K5: 4.9 -> 4.8

# Row 6
$ws.Cells.Item(6, 9).Value  = 1.49   # I6: 1.5 -> 1.49
$ws.Cells.Item(6, 11).Value = 5.2    # K6: 5.4 -> 5.2
$ws.Cells.Item(6, 18).Value = 1.58   # R6: 1.57 -> 1.58
$ws.Cells.Item(6, 29).Value = 11.5   # AC6: 12 -> 11.5
$ws.Cells.Item(6, 30).Value = 10.5   # AD6: 11 -> 10.5
$ws.Cells.Item(6, 31).Value = 15     # AE6: 15.5 -> 15
$ws.Cells.Item(6, 32).Value = 60     # AF6: 65 -> 60
$ws.Cells.Item(6, 34).Value = 26     # AH6: 980 -> 26
$ws.Cells.Item(6, 37).Value = 120    # AK6: 110 -> 120
$ws.Cells.Item(6, 38).Value = 100    # AL6: 80 -> 100
$ws.Cells.Item(6, 39).Value = 120    # AM6: 130 -> 120

$wb.Save()
